# Pilot_Input_Crop_Calendar.xlsx - "updating part 1 closing issue 8"
#
# The small lookup table in row 1-2 of "First Sheet" (A1:K3) described a
# bimodal crop calendar (two init windows: init1_* and init2_*, plus
# dev/mid/late). It is collapsed to a unimodal calendar: the init1 window
# (columns B:C, "01/11".."31/12") is dropped, the former init2 window
# (columns D:E) becomes the new init window and is relabelled
# init_start/init_end, the dev/mid columns shift left unchanged, and the
# mid_end / late_start / late_end dates move out (30/09->30/10,
# 01/10->01/11, 31/10->30/12). The now-blank helper row 3 is removed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "init1_start"/"init1_end" columns (B1:C2) and shift the
# remaining init2/dev/mid/late columns left to take their place.
$ws.Range("B1:C2").Delete(-4159)  # xlShiftToLeft

# The columns that used to be init2_start/init2_end (now B:C) keep their
# date values but get the collapsed-to-unimodal header labels.
$ws.Range("B1").Value = "init_start"
$ws.Range("C1").Value = "init_end"

# mid_end / late_start / late_end dates shift later.
$ws.Range("G2").Value = "30/10"
$ws.Range("H2").Value = "01/11"
$ws.Range("I2").Value = "30/12"

# The leftover blank formatted row is removed.
$ws.Rows("3:3").Delete()

# Match the author's final selection.
$ws.Range("I2").Select()
